$d = $word.ActiveDocument

$tbl = $d.Tables.Item(1)

$map = @{
    "DCL01" = "Clientes"
    "DUS01" = "Usuarios"
    "DPR01" = "Productos"
    "DSE01" = "Servicios"
    "DCR01" = "Credenciales"
    "DPF01" = "Personas Físicas"
    "DPJ01" = "Personas Jurídicas"
}

for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $cell1 = $tbl.Cell($i, 1)
    $key = $cell1.Range.Text.Trim()
    $key = $key -replace "[\x07\x0d\x0a]", ""
    if ($map.ContainsKey($key)) {
        $d.Tables.Item(1).Cell($i, 2).Range.Text = $map[$key]
        $d.Tables.Item(1).Cell($i, 2).Range.LanguageID = 3082
    }
}
